$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")
$ws.Activate()

# Expand table Tabla1 to include row 3 (as when typing into the row right below the table)
$lo = $ws.ListObjects.Item("Tabla1")
$lo.Resize($ws.Range("A1:D3"))

# Fill row 3 (part of the table) - name/surname/email first
$ws.Range("A3").Value = "Fabian"
$ws.Range("B3").Value = "Guia"
$ws.Range("C3").Value = "prueba@gmail.com"

# Fill row 4 (outside the table, plain cells)
$ws.Range("A4").Value = "Aram"
$ws.Range("B4").Value = "Gonzales"
$ws.Range("C4").Value = "correogenerico@gmail.com"
$ws.Range("D4").Value = "algo1234"

# Password placeholder for row 3 typed last
$ws.Range("D3").Value = "contraseña"

# Update selection to G6
$ws.Range("G6").Select() | Out-Null
